$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing cell styles, force Text format so numeric-looking
# strings (e.g. "212.94", "1.00") are not coerced into Number cells,
# then restore the original style/number-format afterwards.
$origStyle = $ws.Range("D2:E51").Style
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = "26.163.72"
$ws.Range("E2").Value2 = "  +3.62%  "
$ws.Range("D3").Value2 = "1.601.91"
$ws.Range("E3").Value2 = "  +2.90%  "
$ws.Range("E4").Value2 = "  -0.32%  "
$ws.Range("D5").Value2 = "212.94"
$ws.Range("E5").Value2 = "  +3.29%  "
$ws.Range("D6").Value2 = "1.00"
$ws.Range("E6").Value2 = "  -0.44%  "
$ws.Range("E7").Value2 = "  +1.98%  "
$ws.Range("E8").Value2 = "  +3.84%  "
$ws.Range("D9").Value2 = "0.0616"
$ws.Range("E9").Value2 = "  +2.33%  "
$ws.Range("E10").Value2 = "  +1.87%  "
$ws.Range("D11").Value2 = "0.0818"
$ws.Range("E11").Value2 = "  +4.83%  "
$ws.Range("D12").Value2 = "1.822.02"
$ws.Range("E12").Value2 = "  +2.75%  "
$ws.Range("D13").Value2 = "1.599.70"
$ws.Range("E13").Value2 = "  +2.68%  "
$ws.Range("D14").Value2 = "4.00"
$ws.Range("E14").Value2 = "  +0.72%  "
$ws.Range("E15").Value2 = "  +2.40%  "
$ws.Range("D16").Value2 = "26.158.51"
$ws.Range("E16").Value2 = "  +3.80%  "
$ws.Range("D17").Value2 = "60.55"
$ws.Range("E17").Value2 = "  +3.03%  "
$ws.Range("D18").Value2 = "0.0₃0721"
$ws.Range("E18").Value2 = "  +2.33%  "
$ws.Range("D20").Value2 = "205.13"
$ws.Range("E20").Value2 = "  +11.18%  "
$ws.Range("D21").Value2 = "4.24"
$ws.Range("E21").Value2 = "  +3.87%  "
$ws.Range("D22").Value2 = "9.31"
$ws.Range("E22").Value2 = "  +1.12%  "
$ws.Range("E23").Value2 = "  +2.75%  "
$ws.Range("E24").Value2 = "  +10.79%  "
$ws.Range("D25").Value2 = "141.69"
$ws.Range("E25").Value2 = "  +1.82%  "
$ws.Range("E26").Value2 = "  -0.44%  "
$ws.Range("E27").Value2 = "  -2.08%  "
$ws.Range("E28").Value2 = "  +3.35%  "
$ws.Range("D29").Value2 = "6.45"
$ws.Range("E29").Value2 = "  +0.81%  "
$ws.Range("E30").Value2 = "  +1.87%  "
$ws.Range("E31").Value2 = "  +2.31%  "
$ws.Range("E32").Value2 = "  +4.10%  "
$ws.Range("D33").Value2 = "2.97"
$ws.Range("E33").Value2 = "  +0.14%  "
$ws.Range("E34").Value2 = "  +2.37%  "
$ws.Range("D35").Value2 = "2.34"
$ws.Range("E35").Value2 = "  +1.63%  "
$ws.Range("D36").Value2 = "1.111.68"
$ws.Range("E36").Value2 = "  +2.58%  "
$ws.Range("D37").Value2 = "0.0162"
$ws.Range("E37").Value2 = "  +9.27%  "
$ws.Range("E38").Value2 = "  -0.05%  "
$ws.Range("D39").Value2 = "2.33"
$ws.Range("E39").Value2 = "  +1.53%  "
$ws.Range("E40").Value2 = "  +2.86%  "
$ws.Range("E41").Value2 = "  +0.85%  "
$ws.Range("E42").Value2 = "  -4.00%  "
$ws.Range("D43").Value2 = "1.735.36"
$ws.Range("E43").Value2 = "  +2.75%  "
$ws.Range("D44").Value2 = "92.92"
$ws.Range("E44").Value2 = "  +0.62%  "
$ws.Range("E45").Value2 = "  +1.46%  "
$ws.Range("E46").Value2 = "  +5.76%  "
$ws.Range("D47").Value2 = "53.41"
$ws.Range("E47").Value2 = "  +2.46%  "
$ws.Range("D48").Value2 = "0.0504"
$ws.Range("E48").Value2 = "  -0.01%  "
$ws.Range("E49").Value2 = "  +0.91%  "
$ws.Range("E50").Value2 = "  -0.29%  "
$ws.Range("E51").Value2 = "  +2.46%  "

$ws.Range("D2:E51").Style = $origStyle
